$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.567.99"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.634.52"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "324.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.526"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "3.048.83"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "2.633.24"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "49.470.60"
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.57%  "
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.137"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("E44").Value = "  -4.29%  "
$ws.Range("D45").Value = "2.059.31"
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.56%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.48%  "
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.72%  "
